# Add a "Save" column (column H) to the s_vals sheet, matching the
# formatting of the existing header cells and filling in the values for
# each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the neighboring header cell (G1, style "s=1":
# bold, bordered, centered) onto the new header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for each data row (unstyled, like the other
# numeric columns).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
